# Insert a new weekly price record for "Cebollín" (Terminal La Palmera de La
# Serena) as row 151, pushing the existing rows 151-234 down to 152-235.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 151:234 down by inserting a new blank row at 151.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row with the latest week's data.
$ws.Range("A151").Value = 8
$ws.Range("B151").Value = "Terminal La Palmera de La Serena"
$ws.Range("C151").Value = "Coquimbo"
$ws.Range("D151").Value = 44806
$ws.Range("E151").Value = 4
$ws.Range("F151").Value = 100112037
$ws.Range("G151").Value = "Cebollín"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 1600
$ws.Range("K151").Value = 1400
$ws.Range("L151").Value = 1600
$ws.Range("M151").Value = 1500
$ws.Range("N151").Value = "`$/paquete 6 unidades"
$ws.Range("O151").Value = "Provincia del Elquí"
$ws.Range("P151").Value = 250
$ws.Range("Q151").Value = 6
$ws.Range("R151").Value = "Hortaliza"
